$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.867.91'
$ws.Range("E2").Value = '  -3.50%  '
$ws.Range("D3").Value = '3.829.79'
$ws.Range("E3").Value = '  -3.21%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''599.77'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = '''166.53'
$ws.Range("E6").Value = '  -2.62%  '
$ws.Range("D7").Value = '3.826.17'
$ws.Range("E7").Value = '  -3.19%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("E10").Value = '  -4.58%  '
$ws.Range("D11").Value = '''6.44'
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("E12").Value = '  -3.47%  '
$ws.Range("D13").Value = '''0.0000259'
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = '''36.80'
$ws.Range("E14").Value = '  -4.88%  '
$ws.Range("D15").Value = '4.471.23'
$ws.Range("E15").Value = '  -3.17%  '
$ws.Range("D16").Value = '3.831.09'
$ws.Range("E16").Value = '  -3.96%  '
$ws.Range("D17").Value = '68.002.10'
$ws.Range("E17").Value = '  -3.17%  '
$ws.Range("D18").Value = '''18.12'
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").Value = '''7.35'
$ws.Range("E19").Value = '  -4.56%  '
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("D22").Value = '''464.57'
$ws.Range("E22").Value = '  -6.86%  '
$ws.Range("D23").Value = '''0.728'
$ws.Range("E23").Value = '  -2.61%  '
$ws.Range("E24").Value = '  -4.95%  '
$ws.Range("D25").Value = '''82.81'
$ws.Range("E25").Value = '  -4.04%  '
$ws.Range("D26").Value = '''2.22'
$ws.Range("E26").Value = '  -4.70%  '
$ws.Range("D27").Value = '''12.02'
$ws.Range("E27").Value = '  -4.03%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '''0.997'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '''9.96'
$ws.Range("E29").Value = '  -3.67%  '
$ws.Range("D30").Value = '''2.95'
$ws.Range("E30").Value = '  -2.27%  '
$ws.Range("D31").Value = '3.978.52'
$ws.Range("E31").Value = '  -3.11%  '
$ws.Range("D32").Value = '''7.54'
$ws.Range("D34").Value = '''31.00'
$ws.Range("E34").Value = '  -4.75%  '
$ws.Range("D35").Value = '''9.51'
$ws.Range("E35").Value = '  -1.84%  '
$ws.Range("D36").Value = '3.789.46'
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("E37").Value = '  -4.96%  '
$ws.Range("D38").Value = '''3.56'
$ws.Range("E38").Value = '  +8.37%  '
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").Value = '''1.02'
$ws.Range("E39").Value = '  -2.72%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '''0.140'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("E41").Value = '  -4.85%  '
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").Value = '''0.310'
$ws.Range("E43").Value = '  -6.13%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '''1.97'
$ws.Range("E44").Value = '  -7.41%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = '''418.50'
$ws.Range("E45").Value = '  -5.30%  '
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").Value = '''0.000294'
$ws.Range("E46").Value = '  +5.46%  '
$ws.Range("E48").Value = '  -1.00%  '
$ws.Range("E49").Value = '  -2.81%  '
$ws.Range("D50").Value = '''141.53'
$ws.Range("E50").Value = '  -1.64%  '
$ws.Range("D51").Value = '''25.99'
$ws.Range("E51").Value = '  +0.63%  '
